$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 277.46155
$ws.Range("I28").Value = 133.41667
$ws.Range("K28").Value = 133.41667
$ws.Range("M28").Value = 351.58333
$ws.Range("H64").Value = 3140
$ws.Range("J64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("N64").ClearContents()
$ws.Range("H67").Value = 3140
$ws.Range("J67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("N67").ClearContents()
$ws.Range("H100").Value = 2660.818
$ws.Range("I100").Value = 1752.7142
$ws.Range("J100").Value = 4250
$ws.Range("K100").Value = 1752.7142
$ws.Range("L100").Value = 4250
$ws.Range("M100").Value = -1211.7142
$ws.Range("N100").Value = -5332
$ws.Range("H108").Value = 24800
$ws.Range("J108").Value = 24800
$ws.Range("L108").Value = 24800
$ws.Range("N108").Value = -32480
$ws.Range("H113").Value = 100004696
$ws.Range("J113").Value = 7996
$ws.Range("L113").Value = 7996
$ws.Range("N113").Value = -14504
$ws.Range("H129").Value = 1842.96
$ws.Range("J129").Value = 2067.1365
$ws.Range("L129").Value = 6201.4095
$ws.Range("N129").Value = -16201.4095
$ws.Range("H141").Value = 1530.579
$ws.Range("I141").Value = 1286.2727
$ws.Range("J141").Value = 3143
$ws.Range("K141").Value = 3858.8181
$ws.Range("L141").Value = 9429
$ws.Range("M141").Value = 1321.1819
$ws.Range("N141").Value = -19789

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1035.6086
$ws.Range("I2").Value = 1012.8823
$ws.Range("K2").Value = 1012.8823
$ws.Range("M2").Value = -899.8823
$ws.Range("H45").Value = 2538.0881
$ws.Range("I45").Value = 2240.2917
$ws.Range("K45").Value = 2240.2917
$ws.Range("M45").Value = -1863.2917
$ws.Range("H61").Value = 2965.6667
$ws.Range("I61").Value = 2469.353
$ws.Range("J61").Value = 4171
$ws.Range("K61").Value = 2469.353
$ws.Range("L61").Value = 4171
$ws.Range("M61").Value = -2257.353
$ws.Range("N61").Value = -4595
$ws.Range("H116").Value = 1035.6086
$ws.Range("I116").Value = 1012.8823
$ws.Range("K116").Value = 1012.8823
$ws.Range("M116").Value = 1281.1177
$ws.Range("H122").Value = 1898.6666
$ws.Range("I122").Value = 1811.5
$ws.Range("J122").Value = 2334.5
$ws.Range("K122").Value = 5434.5
$ws.Range("L122").Value = 7003.5
$ws.Range("M122").Value = -2984.5
$ws.Range("N122").Value = -11903.5
$ws.Range("H132").Value = 15069.3
$ws.Range("I132").Value = 2190.6
$ws.Range("J132").Value = 53705.4
$ws.Range("K132").Value = 6571.799999999999
$ws.Range("L132").Value = 161116.2
$ws.Range("M132").Value = -4041.799999999999
$ws.Range("N132").Value = -166176.2
$ws.Range("H136").Value = 2965.6667
$ws.Range("I136").Value = 2469.353
$ws.Range("J136").Value = 4171
$ws.Range("K136").Value = 7408.059
$ws.Range("L136").Value = 12513
$ws.Range("M136").Value = -4858.059
$ws.Range("N136").Value = -17613

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1035.6086
$ws.Range("I3").Value = 1012.8823
$ws.Range("K3").Value = 1012.8823
$ws.Range("M3").Value = -898.8823
$ws.Range("H134").Value = 3519.432
$ws.Range("I134").Value = 3711.375
$ws.Range("J134").Value = 1600
$ws.Range("K134").Value = 11134.125
$ws.Range("L134").Value = 4800
$ws.Range("M134").Value = -8599.125
$ws.Range("N134").Value = -9870

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1312.4445
$ws.Range("I16").Value = 1117.4286
$ws.Range("K16").Value = 1117.4286
$ws.Range("M16").Value = -830.4286
$ws.Range("H31").Value = 3637.362
$ws.Range("I31").Value = 1816.5454
$ws.Range("J31").Value = 6040.84
$ws.Range("K31").Value = 1816.5454
$ws.Range("L31").Value = 6040.84
$ws.Range("M31").Value = -1521.5454
$ws.Range("N31").Value = -6630.84
$ws.Range("H34").Value = 3637.362
$ws.Range("I34").Value = 1816.5454
$ws.Range("J34").Value = 6040.84
$ws.Range("K34").Value = 1816.5454
$ws.Range("L34").Value = 6040.84
$ws.Range("M34").Value = -1614.5454
$ws.Range("N34").Value = -6444.84
$ws.Range("H52").Value = 21949.908
$ws.Range("J52").Value = 25137.777
$ws.Range("L52").Value = 25137.777
$ws.Range("N52").Value = -25725.777
$ws.Range("H58").Value = 17116.344
$ws.Range("I58").Value = 1546.9333
$ws.Range("J58").Value = 30854.059
$ws.Range("K58").Value = 1546.9333
$ws.Range("L58").Value = 30854.059
$ws.Range("M58").Value = -1343.9333
$ws.Range("N58").Value = -31260.059
$ws.Range("H99").Value = 3355.8823
$ws.Range("I99").Value = 2625
$ws.Range("J99").Value = 4178.125
$ws.Range("K99").Value = 2625
$ws.Range("L99").Value = 4178.125
$ws.Range("M99").Value = -1127
$ws.Range("N99").Value = -7174.125
$ws.Range("H113").Value = 1312.4445
$ws.Range("I113").Value = 1117.4286
$ws.Range("K113").Value = 1117.4286
$ws.Range("M113").Value = 1052.5714
$ws.Range("H122").Value = 995.7241
$ws.Range("I122").Value = 836.94446
$ws.Range("K122").Value = 2510.83338
$ws.Range("M122").Value = -60.83338000000003
$ws.Range("H126").Value = 3355.8823
$ws.Range("I126").Value = 2625
$ws.Range("J126").Value = 4178.125
$ws.Range("K126").Value = 7875
$ws.Range("L126").Value = 12534.375
$ws.Range("M126").Value = -5405
$ws.Range("N126").Value = -17474.375
$ws.Range("H132").Value = 3798.389
$ws.Range("I132").Value = 2920.3635
$ws.Range("J132").Value = 5178.143
$ws.Range("K132").Value = 8761.0905
$ws.Range("L132").Value = 15534.429
$ws.Range("M132").Value = -6231.0905
$ws.Range("N132").Value = -20594.429
$ws.Range("H134").Value = 1313.1177
$ws.Range("I134").Value = 1054.8
$ws.Range("K134").Value = 3164.4
$ws.Range("M134").Value = -629.3999999999996
$ws.Range("H136").Value = 17116.344
$ws.Range("I136").Value = 1546.9333
$ws.Range("J136").Value = 30854.059
$ws.Range("K136").Value = 4640.7999
$ws.Range("L136").Value = 92562.177
$ws.Range("M136").Value = -2090.7999
$ws.Range("N136").Value = -97662.177

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1246.5938
$ws.Range("J5").Value = 2083.111
$ws.Range("L5").Value = 6249.333
$ws.Range("N5").Value = -6473.333
$ws.Range("H131").Value = 734.78
$ws.Range("J131").Value = 769.2717
$ws.Range("L131").Value = 2307.8151
$ws.Range("N131").Value = -12387.8151
$ws.Range("H135").Value = 1246.5938
$ws.Range("J135").Value = 2083.111
$ws.Range("L135").Value = 18747.999
$ws.Range("N135").Value = -23817.999

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 3925
$ws.Range("I97").Value = 3566.6667
$ws.Range("K97").Value = 3566.6667
$ws.Range("M97").Value = -3070.6667
$ws.Range("H122").Value = 5849.5
$ws.Range("I122").Value = 5832.5
$ws.Range("J122").Value = 5875
$ws.Range("K122").Value = 17497.5
$ws.Range("L122").Value = 17625
$ws.Range("M122").Value = -15047.5
$ws.Range("N122").Value = -22525
$ws.Range("H132").Value = 82262.125
$ws.Range("I132").Value = 99536.95
$ws.Range("J132").Value = 49282.91
$ws.Range("K132").Value = 298610.85
$ws.Range("L132").Value = 147848.73
$ws.Range("M132").Value = -296080.85
$ws.Range("N132").Value = -152908.73

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3200
$ws.Range("I7").Value = 2833.3333
$ws.Range("J7").Value = 4850
$ws.Range("K7").Value = 2833.3333
$ws.Range("L7").Value = 4850
$ws.Range("M7").Value = -2721.3333
$ws.Range("N7").Value = -5074
$ws.Range("H40").Value = 4274.25
$ws.Range("I40").Value = 3844.6365
$ws.Range("K40").Value = 3844.6365
$ws.Range("M40").Value = -3708.6365
$ws.Range("H61").Value = 3524.2
$ws.Range("I61").Value = 2033.8572
$ws.Range("J61").Value = 7001.6665
$ws.Range("K61").Value = 2033.8572
$ws.Range("L61").Value = 7001.6665
$ws.Range("M61").Value = -1831.8572
$ws.Range("N61").Value = -7405.6665
$ws.Range("H100").Value = 2446.077
$ws.Range("I100").Value = 1834
$ws.Range("K100").Value = 1834
$ws.Range("M100").Value = -1293
$ws.Range("H113").Value = 3524.2
$ws.Range("I113").Value = 2033.8572
$ws.Range("J113").Value = 7001.6665
$ws.Range("K113").Value = 2033.8572
$ws.Range("L113").Value = 7001.6665
$ws.Range("M113").Value = 136.1428000000001
$ws.Range("N113").Value = -11341.6665
$ws.Range("H126").Value = 3200
$ws.Range("I126").Value = 2833.3333
$ws.Range("J126").Value = 4850
$ws.Range("K126").Value = 8499.999899999999
$ws.Range("L126").Value = 14550
$ws.Range("M126").Value = -6029.999899999999
$ws.Range("N126").Value = -19490
$ws.Range("H132").Value = 229910.2
$ws.Range("I132").Value = 336904.06
$ws.Range("J132").Value = 3335
$ws.Range("K132").Value = 1010712.18
$ws.Range("L132").Value = 10005
$ws.Range("M132").Value = -1008182.18
$ws.Range("N132").Value = -15065

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H64").Value = 10836.723
$ws.Range("J64").Value = 15024.167
$ws.Range("L64").Value = 15024.167
$ws.Range("N64").Value = -15520.167
$ws.Range("H67").Value = 10836.723
$ws.Range("J67").Value = 15024.167
$ws.Range("L67").Value = 15024.167
$ws.Range("N67").Value = -16740.167
$ws.Range("H109").Value = 25763.5
$ws.Range("J109").Value = 25763.5
$ws.Range("L109").Value = 25763.5
$ws.Range("N109").Value = -28537.5
$ws.Range("H113").Value = 822.7037
$ws.Range("I113").Value = 1098.0667
$ws.Range("K113").Value = 3294.2001
$ws.Range("M113").Value = -1124.2001
$ws.Range("H126").Value = 1920.5667
$ws.Range("I126").Value = 1532.1428
$ws.Range("K126").Value = 4596.428400000001
$ws.Range("M126").Value = -2126.428400000001
